$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("FrameCounts")

# Row 45
$ws1.Range("A45").Value = "Begin walljump"
$ws1.Range("B45").Value = 5294
$ws1.Range("C45").Value = 5492

# Row 46
$ws1.Range("A46").Value = "X = 188"
$ws1.Range("B46").Value = 5475
$ws1.Range("C46").Value = 5682

# Row 47
$ws1.Range("A47").Value = "X = 822"
$ws1.Range("B47").Value = 5912
$ws1.Range("C47").Value = 6210

# Row 48
$ws1.Range("A48").Value = "X = 1057"
$ws1.Range("B48").Value = 5967
$ws1.Range("C48").Value = 6265

# Row 49
$ws1.Range("A49").Value = "Black screen"
$ws1.Range("B49").Value = 6162
$ws1.Range("C49").Value = 6474

# Row 50
$ws1.Range("A50").Value = "Batman Appears"
$ws1.Range("B50").Value = 6164
$ws1.Range("C50").Value = 6476

# Update the frozen pane / selection on FrameCounts sheet
$ws1.Activate()
$ws1.Range("B51").Select()
